$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Cover Case Bind  4p"
$ws.Range("C2").Value = "Print Uncollated F 4x0"
$ws.Range("B3").Value = "Text 12p - 1"
$ws.Range("B4").Value = "Text 12p - 3"
$ws.Range("C4").Value = "Print Uncollated F/B 4x4"
$ws.Range("B5").Value = "Text  4p"
$ws.Range("B6").Value = "Text 12p - 2"
$ws.Range("B7").Value = "Cover Case Bind  4p"
$ws.Range("C7").Value = "Print Uncollated F 4x0"
$ws.Range("E7").Value = "1.00"
$ws.Range("B8").Value = "Text 12p - 1"
$ws.Range("B9").Value = "Text 12p - 3"
$ws.Range("C9").Value = "Print Uncollated F/B 4x4"
$ws.Range("E9").Value = "2.00"
$ws.Range("B10").Value = "Text  4p"
$ws.Range("B11").Value = "Text 12p - 2"
$ws.Range("B13").Value = "Cover Case Bind  4p"
$ws.Range("C13").Value = "Print Uncollated F 4x0"
$ws.Range("D13").Value = "Magenta - IS29 Inkjet - "
$ws.Range("E13").Value = "0.05"
$ws.Range("B14").Value = "Cover Case Bind  4p"
$ws.Range("C14").Value = "Print Uncollated F 4x0"
$ws.Range("D14").Value = "Yellow - IS29 Inkjet - "
$ws.Range("E14").Value = "0.05"
$ws.Range("B15").Value = "Cover Case Bind  4p"
$ws.Range("C15").Value = "Print Uncollated F 4x0"
$ws.Range("D15").Value = "Black - IS29 Inkjet - "
$ws.Range("E15").Value = "0.05"
$ws.Range("B16").Value = "Cover Case Bind  4p"
$ws.Range("C16").Value = "Print Uncollated F 4x0"
$ws.Range("D16").Value = "Cyan - IS29 Inkjet - "
$ws.Range("E16").Value = "0.05"
$ws.Range("B17").Value = "Text 12p - 1"
$ws.Range("D17").Value = "Yellow - IS29 Inkjet - "
$ws.Range("B18").Value = "Text 12p - 1"
$ws.Range("D18").Value = "Cyan - IS29 Inkjet - "
$ws.Range("B19").Value = "Text 12p - 1"
$ws.Range("D19").Value = "Magenta - IS29 Inkjet - "
$ws.Range("B20").Value = "Text 12p - 1"
$ws.Range("D20").Value = "Black - IS29 Inkjet - "
$ws.Range("B21").Value = "Text 12p - 3"
$ws.Range("C21").Value = "Print Uncollated F/B 4x4"
$ws.Range("D21").Value = "Yellow - IS29 Inkjet - "
$ws.Range("E21").Value = "0.21"
$ws.Range("B22").Value = "Text 12p - 3"
$ws.Range("C22").Value = "Print Uncollated F/B 4x4"
$ws.Range("D22").Value = "Cyan - IS29 Inkjet - "
$ws.Range("E22").Value = "0.21"
$ws.Range("B23").Value = "Text 12p - 3"
$ws.Range("C23").Value = "Print Uncollated F/B 4x4"
$ws.Range("D23").Value = "Magenta - IS29 Inkjet - "
$ws.Range("E23").Value = "0.21"
$ws.Range("B24").Value = "Text 12p - 3"
$ws.Range("C24").Value = "Print Uncollated F/B 4x4"
$ws.Range("D24").Value = "Black - IS29 Inkjet - "
$ws.Range("E24").Value = "0.21"
$ws.Range("B25").Value = "Text  4p"
$ws.Range("D25").Value = "Cyan - IS29 Inkjet - "
$ws.Range("E25").Value = "0.09"
$ws.Range("B26").Value = "Text  4p"
$ws.Range("D26").Value = "Yellow - IS29 Inkjet - "
$ws.Range("E26").Value = "0.09"
$ws.Range("B27").Value = "Text  4p"
$ws.Range("D27").Value = "Black - IS29 Inkjet - "
$ws.Range("E27").Value = "0.09"
$ws.Range("B28").Value = "Text  4p"
$ws.Range("D28").Value = "Magenta - IS29 Inkjet - "
$ws.Range("E28").Value = "0.09"
$ws.Range("B29").Value = "Text 12p - 2"
$ws.Range("D29").Value = "Yellow - IS29 Inkjet - "
$ws.Range("B30").Value = "Text 12p - 2"
$ws.Range("D30").Value = "Cyan - IS29 Inkjet - "
$ws.Range("B31").Value = "Text 12p - 2"
$ws.Range("D31").Value = "Magenta - IS29 Inkjet - "
$ws.Range("B32").Value = "Text 12p - 2"
$ws.Range("D32").Value = "Black - IS29 Inkjet - "
$ws.Range("B33").Value = "End Sheet  4p - 1"
$ws.Range("C33").Value = "Cut for Press"
$ws.Range("D33").Value = "Special Order Uncoated Cover Non FSC 80# 26 x 40`" 210 ppi"
$ws.Range("E33").Value = "145.00"
$ws.Range("F33").Value = "Sht."
$ws.Range("H33").Value = "429 - Offset Stock cost ea for Jobs Only-Use  for shipping too."
$ws.Range("B34").Value = "Text 12p - 1"
$ws.Range("D34").Value = "Cougar Opaque Smooth Text White Domtar FSC 100# 28 x 40`" 260 ppi"
$ws.Range("E34").Value = "582.00"
$ws.Range("B35").Value = "Text 12p - 3"
$ws.Range("D35").Value = "Cougar Opaque Smooth Text White Domtar FSC 100# 28 x 40`" 260 ppi"
$ws.Range("B36").Value = "Text  4p"
$ws.Range("D36").Value = "Cougar Opaque Smooth Text White Domtar FSC 100# 28 x 40`" 260 ppi"
$ws.Range("E36").Value = "189.00"
$ws.Range("F36").Value = "Sht"
$ws.Range("H36").Value = "430 - Offset Stock Per 1000 cost for Jobs Only"
$ws.Range("B37").Value = "Text 12p - 2"
$ws.Range("D37").Value = "Cougar Opaque Smooth Text White Domtar FSC 100# 28 x 40`" 260 ppi"
$ws.Range("B38").Value = "End Sheet  4p - 2"
$ws.Range("B39").Value = "Cover Case Bind  4p"
$ws.Range("D39").Value = "Sterling Ultra C1S Gloss Verso Non FSC 100# 28 x 40`" 400 ppi"
$ws.Range("E39").Value = "383.00"
